# Applies the DaCapo / JDK21 / ShenandoahGC README stats fix to the single
# one-column results table in the document.
#
# The table holds one value per row. A handful of summary rows (1-3, 4,
# 6-12) get corrected numbers, and three rows further down (44-46) that
# used to carry a whole tab-separated "21 <tab> ... <tab> 100.0" line get
# collapsed back down to the single corrected summary value that row
# originally had (99.97 / 0.02 / 63).

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Set-CellValue($table, $rowIndex, $expectedOld, $newValue) {
    $cell = $table.Cell($rowIndex, 1)
    # Cell.Range.Text carries the trailing end-of-cell marker (CR + BEL);
    # strip trailing control chars before comparing against the plain
    # expected text.
    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)
    if ($current -ne $expectedOld) {
        throw "Row $rowIndex : expected '$expectedOld' but found '$current'"
    }
    $cell.Range.Text = $newValue
}

Set-CellValue $tbl 1  "99.97"   "0M"
Set-CellValue $tbl 2  "0.02"    "0M"
Set-CellValue $tbl 3  "63"      "0M"
Set-CellValue $tbl 4  "21"      "84"

Set-CellValue $tbl 6  "0.00011" "0.00059"
Set-CellValue $tbl 7  "0.00009" "0.00021"
Set-CellValue $tbl 8  "0.00002" "0.00006"
Set-CellValue $tbl 9  "0.00007" "0.00030"
Set-CellValue $tbl 10 "0.00009" "0.00040"
Set-CellValue $tbl 11 "0.00010" "0.00045"
Set-CellValue $tbl 12 "0.00181" "0.01790"

Set-CellValue $tbl 44 "21`t0.00021`t0.00059`t0.00038`t0.00010`t0.00030`t0.00040`t0.00045`t0.00798`t100.0" "99.97"
Set-CellValue $tbl 45 "21`t0.00013`t0.00046`t0.00019`t0.00007`t0.00015`t0.00018`t0.00021`t0.00407`t100.0" "0.02"
Set-CellValue $tbl 46 "21`t0.00013`t0.00039`t0.00019`t0.00006`t0.00015`t0.00019`t0.00022`t0.00404`t100.0" "63"

Write-Output "done"
